# Modelo relacional atualizado. Name adaptado para first_name e last_name
#
# A coluna unica "NOME COMPLETO DO ALUNO" vira duas colunas separadas:
#   G = "NOME DO ALUNO"       (first_name)
#   H = "SOBRENOME DO ALUNO"  (last_name)  <- nova coluna inserida
# As colunas antigas H (DATA DE NASCIMENTO) e I (SEXO) - e os respectivos
# valores da linha 2 - deslizam uma posicao para a direita (I e J).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere uma nova coluna em H, empurrando as colunas seguintes para a
# direita (DATA DE NASCIMENTO: H->I, SEXO: I->J, valor "Menino": I2->J2).
$ws.Columns.Item(8).Insert()

# A nova coluna H deve ficar com a mesma formatacao de cabecalho da coluna
# "DATA DE NASCIMENTO" (agora em I), entao copiamos o formato de I1 para H1.
$ws.Cells.Item(1, 9).Copy()
$ws.Cells.Item(1, 8).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Atualiza os textos dos cabecalhos envolvidos na divisao da coluna.
$ws.Cells.Item(1, 7).Value = "NOME DO ALUNO"
$ws.Cells.Item(1, 8).Value = "SOBRENOME DO ALUNO"

# Largura da nova coluna H, em linha com as colunas vizinhas G/I.
$ws.Columns.Item(8).ColumnWidth = 50

# Reflete a posicao de rolagem/selecao salva na planilha.
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("H6").Select()
